$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so values such as
# "1.000" or "238.48" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.397.55'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.878.04'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '238.48'
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '0.4778'
$ws.Range("E7").Value = '  -1.25%  '
$ws.Range("D8").Value = '0.2826'
$ws.Range("E8").Value = '  -2.62%  '
$ws.Range("D9").Value = '0.06523'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("D10").Value = '1.874.81'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").Value = '0.07466'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '16.67'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").Value = '5.098'
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("D14").Value = '88.27'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").Value = '0.6609'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '30.379.90'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '13.31'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").Value = '0.000007613'
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("D20").Value = '2.113.38'
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("D21").Value = '5.303'
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '218.91'
$ws.Range("E23").Value = '  +13.39%  '
$ws.Range("D24").Value = '6.213'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = '9.376'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").Value = '167.60'
$ws.Range("E26").Value = '  +2.35%  '
$ws.Range("D27").Value = '18.46'
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = '1.977'
$ws.Range("E28").Value = '  +1.73%  '
$ws.Range("D29").Value = '1.463'
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("D30").Value = '0.09383'
$ws.Range("E30").Value = '  +2.78%  '
$ws.Range("D31").Value = '4.322'
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("D32").Value = '4.039'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = '0.05050'
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").Value = '1.204'
$ws.Range("E34").Value = '  +5.11%  '
$ws.Range("D35").Value = '0.7469'
$ws.Range("E35").Value = '  +1.78%  '
$ws.Range("D36").Value = '2.712'
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = '0.01823'
$ws.Range("E37").Value = '  +1.72%  '
$ws.Range("D38").Value = '2.614'
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D39").Value = '2.074'
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").Value = '0.9056'
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("D41").Value = '106.93'
$ws.Range("E41").Value = '  +0.76%  '
$ws.Range("D42").Value = '5.903'
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").Value = '0.4283'
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = '7.423'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").Value = '64.39'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("D47").Value = '0.1281'
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("D48").Value = '1.483'
$ws.Range("E48").Value = '  -3.92%  '
$ws.Range("D49").Value = '8.908'
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("D50").Value = '33.74'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").Value = '0.3904'
$ws.Range("E51").Value = '  +1.05%  '
